# Daily attendance processing - 2025-12-31 08:40:58
# Applies the day's attendance-recording pass to the "Session Analysis Results" sheet:
#  1) Updates the top "Class Statistics" summary (Missing / Pending session counts)
#  2) Normalizes the "Recorded By" text for sessions recorded by both the user and System
#  3) Updates the per-group "Group Statistics" Missing/Pending columns for groups whose
#     next scheduled session (31/12/2025) has now passed without being recorded
#  4) Flips that now-passed session row from "Pending" to "Not Recorded" for each of the
#     six affected groups (B1D1, B1D2, B1E1, B1E2, B1F1, B1F2), matching the formatting
#     already used elsewhere in the sheet for "Not Recorded" rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Class Statistics summary block (K2:L10)
$ws.Range("L7").Value = 9    # Missing Sessions: 3 -> 9
$ws.Range("L8").Value = 114  # Pending Sessions: 120 -> 114

# 2) "Recorded By" column: reorder "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
$recordedByRows = @(8, 9, 10, 12, 14, 15, 17, 34, 35, 36, 38, 40, 41, 43, 60, 61, 62, 64, 66, 67, 69, 86, 87, 88, 90, 92, 93, 95, 112, 113, 114, 116, 118, 119, 121, 138, 139, 140, 142, 144, 145, 147, 164, 167, 170, 191, 194, 197, 218, 221, 224, 245, 248, 251, 272, 275, 278, 299, 302, 305)
foreach ($r in $recordedByRows) {
    $ws.Cells.Item($r, 7).Value = "System, dnasr281@gmail.com"
}

# 3) Group Statistics table (K14:S26) - Missing (P) up by 1, Pending (Q) down by 1
#    for groups B1D1, B1D2, B1E1, B1E2, B1F1, B1F2 (rows 21-26)
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 10

$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 10

$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 10

$ws.Range("P24").Value = 2
$ws.Range("Q24").Value = 10

$ws.Range("P25").Value = 1
$ws.Range("Q25").Value = 10

$ws.Range("P26").Value = 1
$ws.Range("Q26").Value = 10

# 4) Flip the 31/12/2025 session row from "Pending" to "Not Recorded" for the six
#    groups whose session on that date has now lapsed. Copy the formatting used by
#    other "Not Recorded" rows (e.g. row 132) so the fill color/style matches, then
#    update the status text.
$pendingToNotRecordedRows = @(174, 201, 228, 255, 282, 309)
$formatSource = $ws.Range("A132:I132")
foreach ($r in $pendingToNotRecordedRows) {
    $formatSource.Copy()
    $dest = $ws.Range("A" + $r + ":I" + $r)
    $dest.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws.Cells.Item($r, 9).Value = "Not Recorded"
}
